$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# New row 7: CheckID 3 - Restore Running
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Maintenance Tasks Running"
$ws.Range("D7").Value = "Restore Running"
$ws.Range("E7").Value = "http://BrentOzar.com/go/backups"

# New row 8: CheckID 4 - Data File Growing
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "SQL Server Internal Maintenance"
$ws.Range("D8").Value = "Data File Growing"
$ws.Range("E8").Value = "http://BrentOzar.com/go/ifi"

# New row 9: CheckID 5 - Long-Running Query Blocking Others
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Query Problems"
$ws.Range("D9").Value = "Long-Running Query Blocking Others"
$ws.Range("E9").Value = "http://BrentOzar.com/go/blocking"

# Hook up the URL cells as real hyperlinks (matching the existing C2 hyperlink pattern)
$ws.Hyperlinks.Add($ws.Range("E7"), "http://BrentOzar.com/go/backups", "", "", "http://BrentOzar.com/go/backups")
$ws.Hyperlinks.Add($ws.Range("E8"), "http://BrentOzar.com/go/ifi", "", "", "http://BrentOzar.com/go/ifi")
$ws.Hyperlinks.Add($ws.Range("E9"), "http://BrentOzar.com/go/blocking", "", "", "http://BrentOzar.com/go/blocking")

# Restore the "next empty row" selection like Excel would leave after data entry
$ws.Range("A10").Select()
